# Append new Lancers listings scraped at 2025-11-29 18:24:04 JST.
# Rows shift down by one (new row inserted at the top of the data, row 2),
# one more new row is inserted in the middle (new row 5), and one new row
# is appended at the end (new row 9). All existing rows get their
# "取得日時" (fetched-at) timestamp refreshed to the new run's timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item('ランサーズ')

$timestamp = '2025-11-29 18:24:04'
$category = 'システム開発'
$deadline = '期限情報なし'

# Final state of every data row (2-9) after the append, top to bottom.
$rows = @(
    @{ B = '【AI系勉強会】「Google Gravity」開発事例発表者募集!個人開発をプレゼンしませんか?'; D = '1,000 ~ 5,000 円 / 固定'; F = 'https://www.lancers.jp/work/detail/5443957'; G = 360; H = '🔥AI,Ai ◆開発' },
    @{ B = '【急募】フロントエンド開発者募集!React/TypeScriptでのシステム構築'; D = '100,000 円 ~ 200,000 円 / 固定'; F = 'https://www.lancers.jp/work/detail/5443491'; G = 323; H = '🔥React,TypeScript ◆開発' },
    @{ B = '【急募】在庫・販売管理ツールの開発依頼'; D = '500,000 円 ~ 1,000,000 円 / 固定'; F = 'https://www.lancers.jp/work/detail/5443889'; G = 170; H = '◆ツール,開発 ◇管理' },
    @{ B = '【急募】革新的ペット向けECプラットフォーム開発エンジニア募集'; D = '200,000 円 ~ 300,000 円 / 固定'; F = 'https://www.lancers.jp/work/detail/5443928'; G = 68; H = '◆開発' },
    @{ B = 'マンション管理組合のシステム設計構築依頼'; D = '500,000 円 ~ 1,000,000 円 / 固定'; F = 'https://www.lancers.jp/work/detail/5443592'; G = 60; H = '◇管理' },
    @{ B = '【Apache Answer構築】弁護士ドットコムのような専門家Q&Aサイトのサーバー構築・初期設定'; D = '50,000 円 ~ 100,000 円 / 固定'; F = 'https://www.lancers.jp/work/detail/5443617'; G = 38; H = '◇サイト' },
    @{ B = '【急募】Wartalesの武器アイコンとモデルを日本刀に差し替え'; D = '20,000 円 ~ 50,000 円 / 固定'; F = 'https://www.lancers.jp/work/detail/5443568'; G = 13; H = $null },
    @{ B = '地方の補助金に詳しい方募集'; D = '10,000 円 ~ 20,000 円 / 固定'; F = 'https://www.lancers.jp/work/detail/5443921'; G = 10; H = $null }
)

# Drop every existing hyperlink on the sheet; they get rebuilt below in the
# correct row order so the relationship ids line up with F2..F9 again.
$ws.Hyperlinks.Delete()

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $row = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $timestamp
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $category
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $deadline
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G

    if ($null -ne $row.H) {
        $ws.Cells.Item($r, 8).Value = $row.H
    }

    $ws.Hyperlinks.Add($ws.Cells.Item($r, 6), $row.F) | Out-Null
}

"Appended 2025-11-29 18:24 JST run; sheet now has $($rows.Count) data rows."
